$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired data (rows 2..17), columns A (Player), B (Position), C (Team)
$data = @(
    @("Darius Garland",      "PG",    "Cleveland Cavaliers"),
    @("Austin Reaves",       "PG,SG", "Los Angeles Lakers"),
    @("Stephen Curry",       "PG,SG", "Golden State Warriors"),
    @("Tyrese Haliburton",   "PG,SG", "Indiana Pacers"),
    @("Franz Wagner",        "SF,PF", "Orlando Magic"),
    @("Daniel Gafford",      "PF,C",  "Dallas Mavericks"),
    @("Kevin Durant",        "SF,PF", "Phoenix Suns"),
    @("Keegan Murray",       "SF,PF", "Sacramento Kings"),
    @("Jarrett Allen",       "C",     "Cleveland Cavaliers"),
    @("Zach Edey",           "C",     "Memphis Grizzlies"),
    @("Trey Murphy III",     "SF,PF", "New Orleans Pelicans"),
    @("Karl-Anthony Towns",  "PF,C",  "New York Knicks"),
    @("Jalen Duren",         "C",     "Detroit Pistons"),
    @("Tyrese Maxey",        "PG,SG", "Philadelphia 76ers"),
    @("OG Anunoby",          "SF,PF", "New York Knicks"),
    @("Mark Williams",       "C",     "Charlotte Hornets")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# The old table had one extra row (18); remove it now that data only spans to row 17
$ws.Rows.Item(18).Delete()
